$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "numberOfWorkers"
$ws.Range("B6").Value = 2
$ws.Range("B6").NumberFormat = "0"

$ws.Range("B7").Select()
